$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Shallow Tube Wells")
$ws2 = $wb.Worksheets.Item("Deep tube wells")

# Shallow Tube Wells sheet updates
$ws1.Range("E2").Value = 81.641929300000001
$ws1.Range("F2").Value = 28.080494300000002
$ws1.Range("E3").Value = 81.636181899999997
$ws1.Range("F3").Value = 28.2576903
$ws1.Range("E4").Value = 81.6307264
$ws1.Range("F4").Value = 28.206494899999999
$ws1.Range("E5").Value = 81.6300983
$ws1.Range("F5").Value = 28.1873921
$ws1.Range("E6").Value = 81.574748099999994
$ws1.Range("F6").Value = 28.192045499999999
$ws1.Range("E7").Value = 81.658814599999999
$ws1.Range("F7").Value = 28.131810000000002
$ws1.Range("E8").Value = 81.647869299999996
$ws1.Range("F8").Value = 28.092858199999998
$ws1.Range("E9").Value = 81.661625799999996
$ws1.Range("F9").Value = 28.049327399999999
$ws1.Range("E10").Value = 81.654347900000005
$ws1.Range("F10").Value = 28.015118300000001
$ws1.Range("E11").Value = 81.713292300000006
$ws1.Range("F11").Value = 27.983303899999999
$ws1.Range("E12").Value = 81.600283599999997
$ws1.Range("F12").Value = 28.024258700000001
$ws1.Range("E13").Value = 81.589428100000006
$ws1.Range("F13").Value = 28.128715100000001
$ws1.Range("E14").Value = 81.551986200000002
$ws1.Range("F14").Value = 28.146545400000001
$ws1.Range("E15").Value = 81.521693499999998
$ws1.Range("F15").Value = 28.111447900000002
$ws1.Range("E17").Value = 81.645064300000001
$ws1.Range("F17").Value = 28.081222400000001
$ws1.Range("E18").Value = 81.597155099999995
$ws1.Range("F18").Value = 28.106038600000002
$ws1.Range("E19").Value = 81.5933074
$ws1.Range("F19").Value = 28.0780663
$ws1.Range("E20").Value = 81.595679500000003
$ws1.Range("F20").Value = 28.082843700000002
$ws1.Range("E21").Value = 81.586179400000006
$ws1.Range("F21").Value = 28.050744300000002
$ws1.Range("E22").Value = 81.652378600000006
$ws1.Range("F22").Value = 28.0020439
$ws1.Range("E23").Value = 81.694684800000005
$ws1.Range("F23").Value = 28.0285434
$ws1.Range("E24").Value = 81.642070000000004
$ws1.Range("F24").Value = 28.080324099999999
$ws1.Range("E25").Value = 81.547736200000003
$ws1.Range("F25").Value = 28.236479299999999
$ws1.Range("E26").Value = 81.479971000000006
$ws1.Range("F26").Value = 28.188046700000001
$ws1.Range("E27").Value = 81.442083600000004
$ws1.Range("F27").Value = 28.289843000000001
$ws1.Range("E28").Value = 81.383682300000004
$ws1.Range("F28").Value = 28.189743
$ws1.Range("E29").Value = 81.337879799999996
$ws1.Range("F29").Value = 28.169159400000002
$ws1.Range("E30").Value = 81.359306799999999
$ws1.Range("F30").Value = 28.305615700000001
$ws1.Range("E31").Value = 81.429742700000006
$ws1.Range("F31").Value = 28.351235599999999
$ws1.Range("E32").Value = 81.310765599999996
$ws1.Range("F32").Value = 28.250470400000001
$ws1.Range("E33").Value = 81.336589599999996
$ws1.Range("F33").Value = 28.301372000000001
$ws1.Range("E34").Value = 81.359149099999996
$ws1.Range("F34").Value = 28.306523500000001
$ws1.Range("E35").Value = 81.262822499999999
$ws1.Range("F35").Value = 28.279571099999998
$ws1.Range("E36").Value = 81.266078399999998
$ws1.Range("F36").Value = 28.3349604
$ws1.Range("E37").Value = 81.5542272
$ws1.Range("F37").Value = 28.0731407
$ws1.Range("E38").Value = 81.554246300000003
$ws1.Range("F38").Value = 28.073109800000001
$ws1.Range("E39").Value = 81.554121199999997
$ws1.Range("F39").Value = 28.073126999999999
$ws1.Range("E41").Value = 81.496865299999996
$ws1.Range("F41").Value = 28.148684500000002
$ws1.Range("E42").Value = 81.479091499999996
$ws1.Range("F42").Value = 28.146412300000001
$ws1.Range("E43").Value = 81.500847300000004
$ws1.Range("F43").Value = 28.1656409
$ws1.Range("E44").Value = 81.567171200000004
$ws1.Range("F44").Value = 28.2282732
$ws1.Range("E45").Value = 81.584886499999996
$ws1.Range("F45").Value = 28.249013099999999
$ws1.Range("E46").Value = 81.288935100000003
$ws1.Range("F46").Value = 28.238539500000002

# Deep tube wells sheet updates
$ws2.Range("E2").Value = 81.5466409
$ws2.Range("F2").Value = 28.0818215
$ws2.Range("E3").Value = 81.521136299999995
$ws2.Range("F3").Value = 28.096043600000002
$ws2.Range("E4").Value = 81.5165796
$ws2.Range("F4").Value = 28.122328100000001
$ws2.Range("E5").Value = 81.517303499999997
$ws2.Range("F5").Value = 28.113236100000002
$ws2.Range("E6").Value = 81.602454100000003
$ws2.Range("F6").Value = 28.0513081
$ws2.Range("E7").Value = 81.516361599999996
$ws2.Range("F7").Value = 28.085180999999999
$ws2.Range("E8").Value = 81.482591299999996
$ws2.Range("F8").Value = 28.105589899999998
$ws2.Range("E9").Value = 81.502630999999994
$ws2.Range("F9").Value = 28.088340599999999
$ws2.Range("E10").Value = 81.511170800000002
$ws2.Range("F10").Value = 28.082853100000001
